# Weekly update to the "Poroto granado" price sheet.
# A new weekly observation is inserted as row 39 (shifting the existing
# rows 39-47 down to 40-48), and the new row is populated with this
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, pushing rows 39-47 down to 40-48.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly observation.
$ws.Range("A39").Value = 4
$ws.Range("B39").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C39").Value = "Los Lagos"
$ws.Range("D39").Value = 44694
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 100112030
$ws.Range("G39").Value = "Poroto granado"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 60
$ws.Range("K39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 30000
$ws.Range("N39").Value = '$/saco 25 kilos'
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 1200
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
